$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure difference_percentage column keeps literal text (e.g. "23.21%") rather than
# being auto-converted to a numeric percentage value.
$ws.Range("E2:E4").NumberFormat = "@"

# Row 2 - register_clicked_register
$ws.Range("A2").Value = "Differences detected in register: register_clicked_register.png_20240807-234144.png vs register_clicked_register.png_20240808-020447.png"
$ws.Range("B2").Value = "Failure"
$ws.Range("C2").Value = "2024-08-08 02:05:09"
$ws.Range("E2").Value = "23.21%"
$ws.Range("F2").Value = "register_clicked_register.png_20240807-234144.png"
$ws.Range("G2").Value = "register_clicked_register.png_20240808-020447.png"

# Row 3 - register_filled_form
$ws.Range("A3").Value = "Differences detected in register: register_filled_form.png_20240807-234147.png vs register_filled_form.png_20240808-020451.png"
$ws.Range("B3").Value = "Failure"
$ws.Range("C3").Value = "2024-08-08 02:05:09"
$ws.Range("E3").Value = "11.59%"
$ws.Range("F3").Value = "register_filled_form.png_20240807-234147.png"
$ws.Range("G3").Value = "register_filled_form.png_20240808-020451.png"

# Row 4 - register_submitted
$ws.Range("A4").Value = "Differences detected in register: register_submitted.png_20240807-234200.png vs register_submitted.png_20240808-020506.png"
$ws.Range("B4").Value = "Failure"
$ws.Range("C4").Value = "2024-08-08 02:05:09"
$ws.Range("E4").Value = "57.25%"
$ws.Range("F4").Value = "register_submitted.png_20240807-234200.png"
$ws.Range("G4").Value = "register_submitted.png_20240808-020506.png"
